$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Mar 2022" block occupies columns AL:AQ (6 cols), mirroring the
# existing per-year blocks (Shareholder Funds, Debts, Total Revenue, PBT,
# PAT, Cash/Cash Eq). Merge header like the other year blocks.
$ws.Range("AL1:AQ1").Merge()
$ws.Range("AL1").Value = "Mar 2022"

$ws.Range("AN2").Value = "Total"
$ws.Range("AO2").Value = "PBT"
$ws.Range("AP2").Value = "PAT"
$ws.Range("AQ2").Value = "Cash"

$ws.Range("AN3").Value = "Revenue"
$ws.Range("AQ3").Value = "Cash Eq"

$ws.Range("AN4").Value = 24.59
$ws.Range("AO4").Value = "1.04"
$ws.Range("AP4").Value = "0.54"
